$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

# The APS (Annual Population Survey) derived rows (employment / self-employment /
# unemployment / inactivity rate + counts) were pointing at a reporting period
# with no published data ("Jan 2023-Dec 2023"). Fix the period text so it
# reflects the latest period that actually has data.
$newPeriod = "Apr 2023 - Mar 2024"
$ws.Range("B2:B9").Value = $newPeriod

# Restore the user's active selection as captured in the saved workbook.
$ws.Range("C8").Select()
